$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# iPP (row 2) diffusion-coefficient data refreshed for the new thesis draft
$ws.Range("B2").Value = 1.5
$ws.Range("D2").Value = 0.00000000148
$ws.Range("E2").Formula = "=6700"

# Column D widened to fit the refreshed values
$ws.Columns.Item(4).ColumnWidth = 16.140625

# Leave the cursor where work left off
$ws.Range("G10").Select()
